# Add the raw per-amplitude "Ind" readings (rows 11-19) that were still only
# averaged on the "9.81" sheet of memory_grav.xlsx. Row 11 currently holds
# AVERAGE(...) formulas - those get replaced with the literal numbers from the
# raw data set, and eight more individual-run rows (12-19) are appended below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 'protectedDiv(sin(add(read(a0, a1), protectedDiv(a2, 0))), limit(a1, conditional(protectedDiv(write(a0, conditional(write(a0, 0, a1), a1), a2), a1), a2), a1))'
$rowData11 = New-Object 'object[,]' 1,16
$rowData11[0,0] = -1601.43
$rowData11[0,1] = -1762.92
$rowData11[0,2] = -1793.76
$rowData11[0,3] = -1806.8
$rowData11[0,4] = -1829.62
$rowData11[0,5] = -1783.36
$rowData11[0,6] = -1754.89
$rowData11[0,7] = -1654.82
$rowData11[0,8] = -1097.29
$rowData11[0,9] = -1537.21
$rowData11[0,10] = -1702.74
$rowData11[0,11] = -1809.08
$rowData11[0,12] = -1903.96
$rowData11[0,13] = -1976.55
$rowData11[0,14] = -2039.45
$rowData11[0,15] = -2085
$ws.Range("B11:Q11").Value = $rowData11

$ws.Range("A12").Value = 'sub(a2, add(sub(sub(sin(a2), read(a0, a2)), read(a0, a1)), protectedDiv(cos(protectedDiv(read(a0, 0), a2)), write(a0, a1, a2))))'
$rowData12 = New-Object 'object[,]' 1,16
$rowData12[0,0] = -777.3
$rowData12[0,1] = -914.73
$rowData12[0,2] = -822.22
$rowData12[0,3] = -979.89
$rowData12[0,4] = -872.65
$rowData12[0,5] = -1091.69
$rowData12[0,6] = -966.03
$rowData12[0,7] = -715.85
$rowData12[0,8] = -485.38
$rowData12[0,9] = -1391.83
$rowData12[0,10] = -1654.53
$rowData12[0,11] = -1883.92
$rowData12[0,12] = -1719.65
$rowData12[0,13] = -1885.24
$rowData12[0,14] = -2013.94
$rowData12[0,15] = -2022.45
$ws.Range("B12:Q12").Value = $rowData12

$ws.Range("A13").Value = 'protectedDiv(limit(sin(0), conditional(protectedDiv(sin(a2), protectedDiv(a2, a1)), read(a0, 0)), sub(conditional(0, limit(0, limit(a2, 0, a2), read(a0, 0))), abs(add(protectedDiv(cos(add(a1, protectedLog(a2))), conditional(a1, a1)), write(a0, 0, a1))))), sin(a2))'
$rowData13 = New-Object 'object[,]' 1,16
$rowData13[0,0] = -84.64
$rowData13[0,1] = -77.21
$rowData13[0,2] = -97.52
$rowData13[0,3] = -176.59
$rowData13[0,4] = -534.52
$rowData13[0,5] = -394.67
$rowData13[0,6] = -431.43
$rowData13[0,7] = -332.02
$rowData13[0,8] = -472.42
$rowData13[0,9] = -313.9
$rowData13[0,10] = -259.23
$rowData13[0,11] = -359.09
$rowData13[0,12] = -1090.83
$rowData13[0,13] = -1540.24
$rowData13[0,14] = -1724.38
$rowData13[0,15] = -1686.35
$ws.Range("B13:Q13").Value = $rowData13

$ws.Range("A14").Value = 'sub(protectedLog(protectedDiv(protectedLog(protectedDiv(a1, conditional(read(a0, read(a0, a1)), a2))), sub(a1, a2))), protectedDiv(conditional(read(a0, 0), write(a0, limit(add(0, write(a0, a2, cos(a2))), 0, 0), a1)), a2))'
$rowData14 = New-Object 'object[,]' 1,16
$rowData14[0,0] = -735.9
$rowData14[0,1] = -757.48
$rowData14[0,2] = -806.37
$rowData14[0,3] = -574.64
$rowData14[0,4] = -586.11
$rowData14[0,5] = -543.74
$rowData14[0,6] = -490.72
$rowData14[0,7] = -438.39
$rowData14[0,8] = -449.81
$rowData14[0,9] = -428.61
$rowData14[0,10] = -769.46
$rowData14[0,11] = -1307.49
$rowData14[0,12] = -1717.2
$rowData14[0,13] = -1833.18
$rowData14[0,14] = -1922.43
$rowData14[0,15] = -1955.6
$ws.Range("B14:Q14").Value = $rowData14

$ws.Range("A15").Value = 'sub(add(add(write(a0, 0, limit(read(a0, conditional(0, a1)), a1, add(a2, a1))), read(a0, read(a0, a1))), read(a0, protectedDiv(abs(0), conditional(cos(abs(a1)), protectedDiv(a2, 0))))), sub(protectedLog(cos(0)), a1))'
$rowData15 = New-Object 'object[,]' 1,16
$rowData15[0,0] = -1994.27
$rowData15[0,1] = -2000.25
$rowData15[0,2] = -1985.17
$rowData15[0,3] = -1935.46
$rowData15[0,4] = -1899.96
$rowData15[0,5] = -1835.84
$rowData15[0,6] = -1779.28
$rowData15[0,7] = -1693.84
$rowData15[0,8] = -1384.36
$rowData15[0,9] = -1799.1
$rowData15[0,10] = -1873.38
$rowData15[0,11] = -2044.61
$rowData15[0,12] = -2064.51
$rowData15[0,13] = -2163.2
$rowData15[0,14] = -2198.24
$rowData15[0,15] = -2222.87
$ws.Range("B15:Q15").Value = $rowData15

$ws.Range("A16").Value = 'sub(read(a0, sin(0)), add(add(protectedDiv(a2, abs(read(a0, abs(protectedDiv(a2, a2))))), add(write(a0, read(a0, 0), a2), protectedDiv(a2, protectedDiv(a2, read(a0, protectedLog(0)))))), a2))'
$rowData16 = New-Object 'object[,]' 1,16
$rowData16[0,0] = -199.75
$rowData16[0,1] = -121.54
$rowData16[0,2] = -221.3
$rowData16[0,3] = -185.66
$rowData16[0,4] = -602.98
$rowData16[0,5] = -689.44
$rowData16[0,6] = -679.79
$rowData16[0,7] = -870.84
$rowData16[0,8] = -1055.8
$rowData16[0,9] = -1147.07
$rowData16[0,10] = -1278.43
$rowData16[0,11] = -1220.82
$rowData16[0,12] = -1226.27
$rowData16[0,13] = -1281.83
$rowData16[0,14] = -1403.98
$rowData16[0,15] = -1589.37
$ws.Range("B16:Q16").Value = $rowData16

$ws.Range("A17").Value = 'sub(sub(read(a0, sin(conditional(0, sin(a1)))), read(a0, write(a0, protectedLog(conditional(sub(a1, 0), 0)), write(a0, a2, protectedDiv(protectedLog(abs(protectedLog(cos(protectedLog(abs(a1)))))), a2))))), protectedDiv(a1, a2))'
$rowData17 = New-Object 'object[,]' 1,16
$rowData17[0,0] = -649.85
$rowData17[0,1] = -573.73
$rowData17[0,2] = -666.32
$rowData17[0,3] = -718.45
$rowData17[0,4] = -610.43
$rowData17[0,5] = -511.37
$rowData17[0,6] = -394.75
$rowData17[0,7] = -390.88
$rowData17[0,8] = -215.67
$rowData17[0,9] = -256.58
$rowData17[0,10] = -402.68
$rowData17[0,11] = -1338.68
$rowData17[0,12] = -1588.12
$rowData17[0,13] = -1705.61
$rowData17[0,14] = -1866.84
$rowData17[0,15] = -1989.53
$ws.Range("B17:Q17").Value = $rowData17

$ws.Range("A18").Value = 'protectedDiv(sub(a2, read(a0, read(a0, a1))), protectedLog(conditional(protectedDiv(a1, cos(a2)), limit(protectedDiv(sin(sub(a2, a2)), protectedLog(limit(a1, protectedDiv(a1, a2), a2))), a1, protectedLog(write(a0, read(a0, cos(conditional(a2, read(a0, a1)))), write(a0, sin(a1), sin(sin(a2)))))))))'
$rowData18 = New-Object 'object[,]' 1,16
$rowData18[0,0] = -255.37
$rowData18[0,1] = -208.79
$rowData18[0,2] = -211.47
$rowData18[0,3] = -189.68
$rowData18[0,4] = -230.32
$rowData18[0,5] = -223.29
$rowData18[0,6] = -232.27
$rowData18[0,7] = -396.58
$rowData18[0,8] = -314.38
$rowData18[0,9] = -392.85
$rowData18[0,10] = -662.52
$rowData18[0,11] = -1159.08
$rowData18[0,12] = -1354.06
$rowData18[0,13] = -1666.42
$rowData18[0,14] = -1833.84
$rowData18[0,15] = -1925.19
$ws.Range("B18:Q18").Value = $rowData18

$ws.Range("A19").Value = 'sub(write(a0, 0, protectedDiv(read(a0, a1), cos(limit(write(a0, cos(abs(limit(a1, abs(abs(0)), 0))), a1), sub(a2, 0), a1)))), write(a0, add(abs(cos(add(a1, a1))), sub(add(sub(sin(sub(0, a1)), 0), 0), a1)), limit(0, protectedDiv(limit(read(a0, read(a0, a2)), read(a0, sub(a2, a2)), sub(0, sub(abs(write(a0, a1, conditional(protectedDiv(a1, protectedLog(cos(a2))), a2))), a2))), a2), protectedDiv(read(a0, add(a2, a2)), protectedLog(cos(a2))))))'
$rowData19 = New-Object 'object[,]' 1,16
$rowData19[0,0] = -134.61
$rowData19[0,1] = -143.61
$rowData19[0,2] = -129.92
$rowData19[0,3] = -152.12
$rowData19[0,4] = -160.94
$rowData19[0,5] = -165.41
$rowData19[0,6] = -178.39
$rowData19[0,7] = -198.05
$rowData19[0,8] = -245.34
$rowData19[0,9] = -297.76
$rowData19[0,10] = -397.98
$rowData19[0,11] = -1401.08
$rowData19[0,12] = -1545.8
$rowData19[0,13] = -1657.87
$rowData19[0,14] = -1744.01
$rowData19[0,15] = -1940.92
$ws.Range("B19:Q19").Value = $rowData19

# Restore the cursor/selection to the cell the author last clicked (F18).
[void]$ws.Range("F18").Select()

